$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 16 new rows before row 32 (old row32 -> 48, old row33 -> 49)
$ws.Range("A32:H47").EntireRow.Insert()

# 2. Copy formatting (styles) from row 31 down into the newly inserted rows
$ws.Range("A31:H31").Copy()
$ws.Range("A32:H47").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Populate new row values (RFID report rows for Hembra/female cattle)
$ws.Cells.Item(32,1).Value = 45082
$ws.Cells.Item(32,2).Value = 21496
$ws.Cells.Item(32,3).Value = 13
$ws.Cells.Item(32,4).Value = 329
$ws.Cells.Item(32,5).Value = "Hembra"
$ws.Cells.Item(32,6).Value = "Gorda"
$ws.Cells.Item(32,7).Value = 0
$ws.Cells.Item(33,1).Value = 45082
$ws.Cells.Item(33,2).Value = 21638
$ws.Cells.Item(33,3).Value = 10
$ws.Cells.Item(33,4).Value = 339
$ws.Cells.Item(33,5).Value = "Hembra"
$ws.Cells.Item(33,6).Value = "Muy Buena"
$ws.Cells.Item(33,7).Value = 0
$ws.Cells.Item(34,1).Value = 45082
$ws.Cells.Item(34,2).Value = 21753
$ws.Cells.Item(34,3).Value = 10
$ws.Cells.Item(34,4).Value = 361
$ws.Cells.Item(34,5).Value = "Hembra"
$ws.Cells.Item(34,6).Value = "Muy Buena"
$ws.Cells.Item(34,7).Value = 0
$ws.Cells.Item(35,1).Value = 45082
$ws.Cells.Item(35,2).Value = 32161
$ws.Cells.Item(35,3).Value = 9
$ws.Cells.Item(35,4).Value = 341
$ws.Cells.Item(35,5).Value = "Hembra"
$ws.Cells.Item(35,6).Value = "Muy Buena"
$ws.Cells.Item(35,7).Value = 0
$ws.Cells.Item(36,1).Value = 45082
$ws.Cells.Item(36,2).Value = 32450
$ws.Cells.Item(36,3).Value = 8
$ws.Cells.Item(36,4).Value = 335
$ws.Cells.Item(36,5).Value = "Hembra"
$ws.Cells.Item(36,6).Value = "Muy Buena"
$ws.Cells.Item(36,7).Value = 0
$ws.Cells.Item(37,1).Value = 45082
$ws.Cells.Item(37,2).Value = 32634
$ws.Cells.Item(37,3).Value = 5
$ws.Cells.Item(37,4).Value = 359
$ws.Cells.Item(37,5).Value = "Hembra"
$ws.Cells.Item(37,6).Value = "Flaca"
$ws.Cells.Item(37,7).Value = 0
$ws.Cells.Item(38,1).Value = 45082
$ws.Cells.Item(38,2).Value = 32725
$ws.Cells.Item(38,3).Value = 7
$ws.Cells.Item(38,4).Value = 332
$ws.Cells.Item(38,5).Value = "Hembra"
$ws.Cells.Item(38,6).Value = "Muy Buena"
$ws.Cells.Item(38,7).Value = 0
$ws.Cells.Item(39,1).Value = 45082
$ws.Cells.Item(39,2).Value = 32807
$ws.Cells.Item(39,3).Value = 7
$ws.Cells.Item(39,4).Value = 342
$ws.Cells.Item(39,5).Value = "Hembra"
$ws.Cells.Item(39,6).Value = "Buena +"
$ws.Cells.Item(39,7).Value = 0
$ws.Cells.Item(40,1).Value = 45082
$ws.Cells.Item(40,2).Value = 33687
$ws.Cells.Item(40,3).Value = 9
$ws.Cells.Item(40,4).Value = 341
$ws.Cells.Item(40,5).Value = "Hembra"
$ws.Cells.Item(40,6).Value = "Muy Buena"
$ws.Cells.Item(40,7).Value = 0
$ws.Cells.Item(41,1).Value = 45082
$ws.Cells.Item(41,2).Value = 33687
$ws.Cells.Item(41,3).Value = 16
$ws.Cells.Item(41,4).Value = 341
$ws.Cells.Item(41,5).Value = "Hembra"
$ws.Cells.Item(41,6).Value = "Gorda"
$ws.Cells.Item(41,7).Value = 0
$ws.Cells.Item(42,1).Value = 45082
$ws.Cells.Item(42,2).Value = 33689
$ws.Cells.Item(42,3).Value = 4
$ws.Cells.Item(42,4).Value = 340
$ws.Cells.Item(42,5).Value = "Hembra"
$ws.Cells.Item(42,6).Value = "Flaca"
$ws.Cells.Item(42,7).Value = 0
$ws.Cells.Item(43,1).Value = 45082
$ws.Cells.Item(43,2).Value = 34123
$ws.Cells.Item(43,3).Value = 12
$ws.Cells.Item(43,4).Value = 356
$ws.Cells.Item(43,5).Value = "Hembra"
$ws.Cells.Item(43,6).Value = "Apenas Gorda"
$ws.Cells.Item(43,7).Value = 0
$ws.Cells.Item(44,1).Value = 45082
$ws.Cells.Item(44,2).Value = 34404
$ws.Cells.Item(44,3).Value = 10
$ws.Cells.Item(44,4).Value = 365
$ws.Cells.Item(44,5).Value = "Hembra"
$ws.Cells.Item(44,6).Value = "Muy Buena"
$ws.Cells.Item(44,7).Value = 0
$ws.Cells.Item(45,1).Value = 45082
$ws.Cells.Item(45,2).Value = 34468
$ws.Cells.Item(45,3).Value = 11
$ws.Cells.Item(45,4).Value = 352
$ws.Cells.Item(45,5).Value = "Hembra"
$ws.Cells.Item(45,6).Value = "Muy Buena"
$ws.Cells.Item(45,7).Value = 0
$ws.Cells.Item(46,1).Value = 45082
$ws.Cells.Item(46,2).Value = 50483
$ws.Cells.Item(46,3).Value = 8
$ws.Cells.Item(46,4).Value = 327
$ws.Cells.Item(46,5).Value = "Hembra"
$ws.Cells.Item(46,6).Value = "Muy Buena"
$ws.Cells.Item(46,7).Value = 0
$ws.Cells.Item(47,1).Value = 45082
$ws.Cells.Item(47,2).Value = 50490
$ws.Cells.Item(47,3).Value = 13
$ws.Cells.Item(47,4).Value = 343
$ws.Cells.Item(47,5).Value = "Hembra"
$ws.Cells.Item(47,6).Value = "Gorda"
$ws.Cells.Item(47,7).Value = 0

# 4. Row heights
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 30.75
$ws.Rows.Item(6).RowHeight = 30.75
$ws.Rows.Item(7).RowHeight = 30.75
$ws.Rows.Item(8).RowHeight = 15.75
$ws.Rows.Item(9).RowHeight = 15.75
$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 15.75
$ws.Rows.Item(13).RowHeight = 15.75
$ws.Rows.Item(14).RowHeight = 15.75
$ws.Rows.Item(15).RowHeight = 15.75
$ws.Rows.Item(16).RowHeight = 15.75
$ws.Rows.Item(17).RowHeight = 15.75
$ws.Rows.Item(18).RowHeight = 30.75
$ws.Rows.Item(19).RowHeight = 15.75
$ws.Rows.Item(20).RowHeight = 15.75
$ws.Rows.Item(21).RowHeight = 15.75
$ws.Rows.Item(22).RowHeight = 15.75
$ws.Rows.Item(23).RowHeight = 15.75
$ws.Rows.Item(24).RowHeight = 15.75
$ws.Rows.Item(25).RowHeight = 15.75
$ws.Rows.Item(26).RowHeight = 15.75
$ws.Rows.Item(27).RowHeight = 15.75
$ws.Rows.Item(28).RowHeight = 15.75
$ws.Rows.Item(29).RowHeight = 15.75
$ws.Rows.Item(30).RowHeight = 15.75
$ws.Rows.Item(31).RowHeight = 15.75
$ws.Rows.Item(48).RowHeight = 30.75
$ws.Rows.Item(49).RowHeight = 15.75

$ws.Rows.Item(32).RowHeight = 15.75
$ws.Rows.Item(33).RowHeight = 15.75
$ws.Rows.Item(34).RowHeight = 15.75
$ws.Rows.Item(35).RowHeight = 15.75
$ws.Rows.Item(36).RowHeight = 15.75
$ws.Rows.Item(37).RowHeight = 15.75
$ws.Rows.Item(38).RowHeight = 15.75
$ws.Rows.Item(39).RowHeight = 15.75
$ws.Rows.Item(40).RowHeight = 15.75
$ws.Rows.Item(41).RowHeight = 15.75
$ws.Rows.Item(42).RowHeight = 15.75
$ws.Rows.Item(43).RowHeight = 30.75
$ws.Rows.Item(44).RowHeight = 15.75
$ws.Rows.Item(45).RowHeight = 15.75
$ws.Rows.Item(46).RowHeight = 15.75
$ws.Rows.Item(47).RowHeight = 15.75

# 5. Update selection to match target view
$ws.Range("E49").Select()
